# Update the "experiments" log worksheet:
#  - Fill in the missing result table for experiment #6 (row 7, column G)
#  - Highlight the best accuracy achieved so far (row 9, column F) in yellow
#  - Add a new experiment #9 row (row 10) for an uncased-model run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NL = [char]10

# --- Row 7: fill in the previously-empty "Best Accuracy (Val-Excel)" result table ---
$g7 = "accuracy    f1_macro    precision-neg    recall-neg" + $NL +
      "----------  ----------  ---------------  ------------" + $NL +
      "91.1%       91.1%       90.05%           92.39%"
$ws.Range("G7").Value = $g7

# --- Row 9: highlight the current best accuracy figure in yellow ---
$ws.Range("F9").Interior.Color = 65535

# --- Row 10: new experiment entry (#9) using the uncased model ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "bert-base-uncased"
$ws.Range("C10").Value = "lowercase"
$ws.Range("D10").Value = "NIL"

$e10 = "seed = 1234" + $NL +
       "test_size = 0.2" + $NL +
       "MAX_LEN = 512" + $NL +
       "START_LEN = 382" + $NL +
       "END_LEN = 128" + $NL +
       "batch_size = 16" + $NL +
       "epochs = 7" + $NL +
       "ATTENTION_PROB_DROPOUT_PROB=0.2" + $NL +
       "HIDDEN_DROPOUT_PROB=0.2" + $NL +
       "SAVE_PROCESSED = True" + $NL +
       "use_gpu_test = True"
$ws.Range("E10").Value = $e10
$ws.Range("E10").WrapText = $true

$ws.Rows.Item(10).RowHeight = 158.4

# Widen column E so the longer config text in the new row is readable
$ws.Columns.Item(5).ColumnWidth = 38.5

# Mirror the author's final selection, as recorded in the saved file
[void]$ws.Range("E10").Select()
